$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.856.07"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.512.06"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'532.43"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'138.57"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D9").Value = "2.515.23"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'5.44"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "'0.357"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "2.956.71"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'23.06"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "58.816.11"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "2.509.85"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'11.03"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "'322.51"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'62.13"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'0.425"
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").Value = "'6.67"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'163.49"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -6.82%  "
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("D39").Value = "'36.80"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'0.802"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'5.21"
$ws.Range("E42").Value = "  -7.84%  "
$ws.Range("D43").Value = "'279.09"
$ws.Range("E43").Value = "  -6.00%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "'10.88"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").Value = "'0.596"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "'0.0933"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'121.65"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'18.42"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -2.43%  "
